$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing covid_deaths (column C) values for various rows ---
$cUpdates = @{
    973  = 13
    974  = 36
    981  = 25
    992  = 34
    997  = 39
    1001 = 31
    1015 = 41
    1020 = 40
    1036 = 38
    1056 = 3
    1057 = 16
    1059 = 45
    1067 = 5
    1068 = 8
    1073 = 17
    1074 = 37
    1084 = 22
    1085 = 45
    1088 = 14
    1089 = 16
    1090 = 31
}

foreach ($row in $cUpdates.Keys) {
    $ws.Range("C$row").Value = $cUpdates[$row]
}

# --- Update existing rows 1091-1094: new agegrp (column B) and, for some, new covid_deaths (column C) ---
$ws.Range("B1091").Value = "30-39"

$ws.Range("B1092").Value = "40-49"
$ws.Range("C1092").Value = 1

$ws.Range("B1093").Value = "50-59"
$ws.Range("C1093").Value = 5

$ws.Range("B1094").Value = "60-69"
$ws.Range("C1094").Value = 6

# --- Append new rows 1095-1101 ---
$newRows = @(
    @{ Row = 1095; Date = 44181; AgeGrp = "70-79"; Deaths = 9 },
    @{ Row = 1096; Date = 44181; AgeGrp = "80+";   Deaths = 32 },
    @{ Row = 1097; Date = 44182; AgeGrp = "0-19";  Deaths = 1 },
    @{ Row = 1098; Date = 44182; AgeGrp = "50-59"; Deaths = 2 },
    @{ Row = 1099; Date = 44182; AgeGrp = "60-69"; Deaths = 2 },
    @{ Row = 1100; Date = 44182; AgeGrp = "70-79"; Deaths = 5 },
    @{ Row = 1101; Date = 44182; AgeGrp = "80+";   Deaths = 12 }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Range("A$r").Value = $item.Date
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("B$r").Value = $item.AgeGrp
    $ws.Range("C$r").Value = $item.Deaths
}
